$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target columns keep their text formatting so numeric-looking
# strings (e.g. "40.90", "1.00") are not coerced into numbers on assignment.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "70.150.96"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "3.738.78"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "623.34"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "180.65"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("D7").Value = "3.736.75"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("E11").Value = "  -4.62%  "
$ws.Range("E12").Value = "  -2.99%  "
$ws.Range("D13").Value = "40.90"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "4.362.59"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "3.736.44"
$ws.Range("D17").Value = "70.166.00"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "505.78"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("D22").Value = "9.35"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").Value = "2.58"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").Value = "86.71"
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("D26").Value = "11.50"
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("D27").Value = "13.19"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("E28").Value = "  +20.74%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "2.49"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").Value = "6.18"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").Value = "0.137"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "0.338"
$ws.Range("E39").Value = "  -2.43%  "
$ws.Range("D40").Value = "2.10"
$ws.Range("E40").Value = "  -6.64%  "
$ws.Range("D41").Value = "50.38"
$ws.Range("E41").Value = "  -2.37%  "
$ws.Range("D42").Value = "45.90"
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("D43").Value = "430.42"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").Value = "8.74"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").Value = "2.87"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").Value = "3.010.64"
$ws.Range("E46").Value = "  -4.46%  "
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("D48").Value = "27.64"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D50").Value = "137.66"
$ws.Range("D51").Value = "2.52"
$ws.Range("E51").Value = "  +1.67%  "
